# Apply cryptos list update (prices + 1h volume %) per commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "58.238.58"
$ws.Range("E2").Value = "  -0.37%  "
$ws.Range("D3").Value = "3.125.72"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  +0.01%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "529.24"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +1.58%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "142.41"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.07%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.01%  "
$ws.Range("D8").Value = "3.126.69"
$ws.Range("E8").Value = "  +1.48%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.445"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +1.90%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "7.16"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("E11").Value = "  +0.48%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.392"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +2.43%  "
$ws.Range("D13").Value = "3.663.50"
$ws.Range("E13").Value = "  +1.38%  "
$ws.Range("E14").Value = "  +3.40%  "
$ws.Range("E15").Value = "  -3.67%  "
$ws.Range("E16").Value = "  -0.09%  "
$ws.Range("D17").Value = "58.302.51"
$ws.Range("E17").Value = "  -0.30%  "
$ws.Range("D18").Value = "3.137.13"
$ws.Range("E18").Value = "  +1.45%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "6.11"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -0.16%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "12.82"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.36%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "7.97"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -1.42%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "343.00"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +0.43%  "
$ws.Range("E23").Value = "  +0.05%  "
$ws.Range("E24").Value = "  +2.08%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "67.69"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +3.00%  "
$ws.Range("E26").Value = "  -0.35%  "
$ws.Range("E27").Value = "  -0.09%  "
$ws.Range("D28").Value = "0.0₃0931"
$ws.Range("E28").Value = "  +2.10%  "
$ws.Range("E29").Value = "  +0.00%  "
$ws.Range("E30").Value = "  +2.11%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "6.40"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -2.49%  "
$ws.Range("E32").Value = "  +2.22%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "21.12"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +0.84%  "
$ws.Range("E34").Value = "  -1.16%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "158.88"
$c.Style = "Normal"
$ws.Range("E35").Value = "  +2.92%  "
$ws.Range("E36").Value = "  +3.76%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "6.20"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +2.75%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "26.34"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -2.02%  "
$ws.Range("E39").Value = "  -3.66%  "
$ws.Range("E40").Value = "  +12.46%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "0.0666"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -1.64%  "
$ws.Range("E42").Value = "  +4.80%  "
$ws.Range("E43").Value = "  +2.65%  "
$ws.Range("D44").Value = "3.167.26"
$ws.Range("E44").Value = "  +1.32%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "36.66"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.15%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +0.00%  "
$ws.Range("E47").Value = "  +2.95%  "
$ws.Range("D48").Value = "2.258.10"
$ws.Range("E48").Value = "  -0.29%  "
$ws.Range("E49").Value = "  +5.28%  "
$ws.Range("E50").Value = "  +2.16%  "
$ws.Range("E51").Value = "  +0.01%  "
